# Rename the existing sheet to "Eco" and add a new "Normal" sheet right
# after it, matching the target workbook layout.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Eco"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Normal"

# --- Populate "Normal" sheet --------------------------------------------
# The "Normal" sheet reuses the same category labels as "Eco" (just a
# different subset/arrangement), so copy whole A:C rows straight from the
# "Eco" sheet. This keeps styles (and shared-string reuse) identical to
# the source workbook instead of fabricating brand-new style entries.
$rowMap = @{
    1  = 1
    2  = 2
    3  = 5
    4  = 4
    5  = 7
    6  = 8
    7  = 5
    8  = 10
    9  = 11
    10 = 12
    11 = 13
    12 = 14
    13 = 16
    14 = 17
    15 = 18
    16 = 19
    17 = 20
}

foreach ($t in 1..17) {
    $s = $rowMap[$t]
    $srcRange = $ws1.Range("A" + $s + ":C" + $s)
    $dstRange = $ws2.Range("A" + $t + ":C" + $t)
    $srcRange.Copy($dstRange)
}

# Rows 3 and 7 are blank spacer rows on "Normal" - their style donor (Eco
# row 5) has text in it, so clear the copied values back out.
$ws2.Range("A3:C3").ClearContents()
$ws2.Range("A7:C7").ClearContents()

# Row heights: 15.75 everywhere except the section-header rows (8, 10, 13,
# 14) which keep the 15pt default.
foreach ($r in 1..17) {
    if ($r -ne 8 -and $r -ne 10 -and $r -ne 13 -and $r -ne 14) {
        $ws2.Rows.Item($r).RowHeight = 15.75
    }
}

# Column widths for the new sheet.
$ws2.Columns.Item(1).ColumnWidth = 19.42578125
$ws2.Columns.Item(2).ColumnWidth = 26.28515625
$ws2.Columns.Item(3).ColumnWidth = 30.85546875

# --- Tidy up "Eco" sheet view -------------------------------------------
# Selection becomes the whole used range instead of the old E7 cell.
$ws1.Range("A1:C20").Select() | Out-Null

# Selection / active cell on "Normal" sits just below the data. Select it
# last so that "Normal" ends up as the active tab, matching the target.
$ws2.Activate()
$ws2.Range("B21").Select() | Out-Null
